# Add six new model-evaluation rows (4-9) below the existing "Mistral-7b
# (temperature = 0.7)" row, matching the header's B/E/H/K/N/Q (+2 blank
# merge-companion cols each) layout, then merge each triplet and move the
# selection to A9 - mirroring how the existing header/row-3 data is laid out.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply the same center alignment used by the existing rows to the whole
# new block first, so the following merges inherit it instead of resetting
# to the default style.
$ws.Range("A4:S9").HorizontalAlignment = -4108

# row -> [ Model name, BERTScore P, BERTScore R, BERTScore F1, ROUGE-1, ROUGE-2, ROUGE-L ]
$rows = @(
    @(4, "deepseek_r1_distill_llama_70b",      0.7854, 0.8183, 0.8013, 0.1067, 0.0145, 0.0599),
    @(5, "deepseek-r1t2-chimera.jsonl",        0.762,  0.8204, 0.7898, 0.109,  0.0173, 0.0618),
    @(6, "gemini",                             0.7695, 0.8269, 0.797,  0.0857, 0.0184, 0.0153),
    @(7, "llama-4-maverick-17b-128e-instruct", 0.7899, 0.8251, 0.8068, 0.1245, 0.0232, 0.0752),
    @(8, "Mistral-7b (temperature = 0.1)",     0.8035, 0.8249, 0.8137, 0.1488, 0.0262, 0.0901),
    @(9, "phi4",                               0.7731, 0.8151, 0.7932, 0.0428, 0.0099, 0.0279)
)

$dataCols = @("B", "E", "H", "K", "N", "Q")

foreach ($r in $rows) {
    $rowNum = $r[0]
    $ws.Range("A$rowNum").Value = $r[1]

    for ($i = 0; $i -lt $dataCols.Length; $i++) {
        $col = $dataCols[$i]
        $value = $r[2 + $i]
        $ws.Range("$col$rowNum").Value = $value

        $nextCol1 = [char]([int][char]$col + 1)
        $nextCol2 = [char]([int][char]$col + 2)
        $ws.Range("$col$rowNum`:$nextCol2$rowNum").Merge()
    }
}

$ws.Range("A9").Select()
